# The workbook gains one new data row: a new row is inserted immediately
# before the current row 358, shifting the existing rows 358:455 down to
# 359:456 (so the sheet dimension grows from A1:R455 to A1:R456), and the
# newly inserted row 358 is populated with a new "Perejil" price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 358; this shifts every row at/after
# 358 down by one (old row 358 -> 359, ..., old row 455 -> 456).
$ws.Rows.Item(358).Insert()

# Populate the newly inserted row 358 with the new record's data.
$ws.Range("A358").Value  = 9
$ws.Range("B358").Value  = "Vega Central Mapocho de Santiago"
$ws.Range("C358").Value  = "Metropolitana"
$ws.Range("D358").Value2 = 44855
$ws.Range("E358").Value  = 13
$ws.Range("F358").Value  = 100112044
$ws.Range("G358").Value  = "Perejil"
$ws.Range("H358").Value  = "Sin especificar"
$ws.Range("I358").Value  = "Primera"
$ws.Range("J358").Value  = 85
$ws.Range("K358").Value  = 9000
$ws.Range("L358").Value  = 10000
$ws.Range("M358").Value  = 9588
$ws.Range("N358").Value  = '$/docena de atados'
$ws.Range("O358").Value  = "Región Metropolitana"
$ws.Range("P358").Value  = 3196
$ws.Range("Q358").Value  = 3
$ws.Range("R358").Value  = "Hortaliza"
